$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.286.25"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.938.45"
$ws.Range("E3").Value = "  -2.76%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "567.19"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -3.33%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "158.81"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.50%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.517"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").Value = "2.936.15"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("E10").Value = "  -3.98%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.150"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.26%  "
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("E13").Value = "  -2.04%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.28"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "65.274.30"
$ws.Range("D17").Value = "3.427.46"
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").Value = "2.940.35"
$ws.Range("E19").Value = "  -2.61%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.94"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +7.85%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "444.96"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.16%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -2.12%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "82.03"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  -1.96%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "12.10"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.93%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.07"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -6.89%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.02"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.38%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.40"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E31").Value = "  -1.97%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0000102"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.80%  "
$ws.Range("E33").Value = "  -0.17%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.111"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.68%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.14%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.971"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.79%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.71"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("E38").Value = "  +0.42%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "44.64"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("E40").Value = "  -9.88%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.83"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -7.53%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.298"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.120"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  -1.04%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "385.16"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "2.705.74"
$ws.Range("E47").Value = "  -3.50%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "133.22"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.91%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.18"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.96%  "
$ws.Range("E51").Value = "  -0.71%  "
